# "Generate Report for Handoff" - mark the a0fcd79e-...md file as handed off
# (status "In Translation" -> "Ready for handoff") for both the zh-cn and
# de-de locales, refreshing the "Latest Handoff Datetime" / "Latest Handoff
# Date" timestamps, and rolling the Overview sheet's summary up to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- zh-cn detail sheet: row 3 is a0fcd79e-4f80-4699-9ce1-af5cda8c14f2.md ---
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-23 12:20:14"

# --- de-de detail sheet: row 3 is a0fcd79e-4f80-4699-9ce1-af5cda8c14f2.md ---
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-23 12:20:19"

# --- Overview sheet: row 3 rolls up both locale statuses + latest date ---
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-23 12:20:19"
